# Weekly update: a new price-report date (2021-10-22, serial 44491) is
# inserted as a new block of 3 rows (Primera/Segunda/Tercera) right before
# the existing row 464, pushing the rest of the "Betarraga" records down by
# three rows (the former last block, rows 527-529, lands on new rows
# 530-532).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 464..466 (existing rows shift down, formatting of
# the row above carries down onto the new rows - matches column D's date
# style).
$ws.Rows("464:466").Insert()

# Shared constant column values for every data row in this table.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100114014
$categoria = "Betarraga"
$variedad = "Sin especificar"
$unidad = "`$/unidad"
$origen = "Región Metropolitana"
$kgUnidades = 1
$clasificacion = "Hortaliza"

$fecha = 44491

# Row 464 - Primera
$ws.Cells.Item(464, 1).Value = $mercadoId
$ws.Cells.Item(464, 2).Value = $mercado
$ws.Cells.Item(464, 3).Value = $region
$ws.Cells.Item(464, 4).Value = $fecha
$ws.Cells.Item(464, 5).Value = $codreg
$ws.Cells.Item(464, 6).Value = $categoriaId
$ws.Cells.Item(464, 7).Value = $categoria
$ws.Cells.Item(464, 8).Value = $variedad
$ws.Cells.Item(464, 9).Value = "Primera"
$ws.Cells.Item(464, 10).Value = 46000
$ws.Cells.Item(464, 11).Value = 90
$ws.Cells.Item(464, 12).Value = 100
$ws.Cells.Item(464, 13).Value = 95
$ws.Cells.Item(464, 14).Value = $unidad
$ws.Cells.Item(464, 15).Value = $origen
$ws.Cells.Item(464, 16).Value = 95
$ws.Cells.Item(464, 17).Value = $kgUnidades
$ws.Cells.Item(464, 18).Value = $clasificacion

# Row 465 - Segunda
$ws.Cells.Item(465, 1).Value = $mercadoId
$ws.Cells.Item(465, 2).Value = $mercado
$ws.Cells.Item(465, 3).Value = $region
$ws.Cells.Item(465, 4).Value = $fecha
$ws.Cells.Item(465, 5).Value = $codreg
$ws.Cells.Item(465, 6).Value = $categoriaId
$ws.Cells.Item(465, 7).Value = $categoria
$ws.Cells.Item(465, 8).Value = $variedad
$ws.Cells.Item(465, 9).Value = "Segunda"
$ws.Cells.Item(465, 10).Value = 34000
$ws.Cells.Item(465, 11).Value = 80
$ws.Cells.Item(465, 12).Value = 85
$ws.Cells.Item(465, 13).Value = 82
$ws.Cells.Item(465, 14).Value = $unidad
$ws.Cells.Item(465, 15).Value = $origen
$ws.Cells.Item(465, 16).Value = 82
$ws.Cells.Item(465, 17).Value = $kgUnidades
$ws.Cells.Item(465, 18).Value = $clasificacion

# Row 466 - Tercera
$ws.Cells.Item(466, 1).Value = $mercadoId
$ws.Cells.Item(466, 2).Value = $mercado
$ws.Cells.Item(466, 3).Value = $region
$ws.Cells.Item(466, 4).Value = $fecha
$ws.Cells.Item(466, 5).Value = $codreg
$ws.Cells.Item(466, 6).Value = $categoriaId
$ws.Cells.Item(466, 7).Value = $categoria
$ws.Cells.Item(466, 8).Value = $variedad
$ws.Cells.Item(466, 9).Value = "Tercera"
$ws.Cells.Item(466, 10).Value = 12000
$ws.Cells.Item(466, 11).Value = 60
$ws.Cells.Item(466, 12).Value = 60
$ws.Cells.Item(466, 13).Value = 60
$ws.Cells.Item(466, 14).Value = $unidad
$ws.Cells.Item(466, 15).Value = $origen
$ws.Cells.Item(466, 16).Value = 60
$ws.Cells.Item(466, 17).Value = $kgUnidades
$ws.Cells.Item(466, 18).Value = $clasificacion
